# Apply the "assertMatch(text,regex)" + "openFile(filePath)" new-command edit
# to the hidden '#system' lookup sheet, and remove the obsolete 'tn.5250'
# target/category together with its whole data column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) "base" commands (column F): insert "assertMatch(text,regex)" in its
#    alphabetically-correct spot, right after "assertEqual(expected,actual)"
#    (row 10) and before "assertNotContain(text,substring)" (row 11).
#    Shift every following value in the column down by one row first.
# ---------------------------------------------------------------------
for ($r = 44; $r -ge 11; $r--) {
    $ws.Cells.Item($r + 1, 6).Value = $ws.Cells.Item($r, 6).Value2
}
$ws.Cells.Item(11, 6).Value = "assertMatch(text,regex)"

# ---------------------------------------------------------------------
# 2) "external" commands (column J): insert "openFile(filePath)" as the
#    new first entry (row 2), shifting the rest down by one row.
# ---------------------------------------------------------------------
for ($r = 6; $r -ge 2; $r--) {
    $ws.Cells.Item($r + 1, 10).Value = $ws.Cells.Item($r, 10).Value2
}
$ws.Cells.Item(2, 10).Value = "openFile(filePath)"

# ---------------------------------------------------------------------
# 3) "target" categories (column A): drop the retired "tn.5250" category
#    (was row 27), shifting everything below it up by one row and
#    clearing the now-empty last row (33).
# ---------------------------------------------------------------------
for ($r = 27; $r -le 32; $r++) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($r + 1, 1).Value2
}
$ws.Cells.Item(33, 1).ClearContents()

# ---------------------------------------------------------------------
# 4) Remove the whole "tn.5250" data column (old column AA); this shifts
#    web / webalert / webcookie / ws / ws.async / xml one column to the
#    left (AB->AA, AC->AB, AD->AC, AE->AD, AF->AE, AG->AF).
# ---------------------------------------------------------------------
$ws.Columns.Item(27).Delete()

# ---------------------------------------------------------------------
# 5) Re-point the named ranges so they describe the new layout.
# ---------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$F`$2:`$F`$45"
$wb.Names.Item("external").RefersTo = "='#system'!`$J`$2:`$J`$7"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("tn.5250").Delete()
$wb.Names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$151"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AF`$2:`$AF`$27"
